$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.338.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.825.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -2.96%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4273'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3695'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07244'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8633'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.904.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.673'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.340'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07105'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008857'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.371.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.163'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.115.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.013'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.259'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08857'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.78%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.204'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7624'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.492'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.844'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.120'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01959'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05264'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.876'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.073'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1681'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5062'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.618'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4723'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06416'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.666'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.812'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.32%  '
